$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1251282873

$ws.Range("I2").Value = "'03/02/2021"

$ws.Range("X2").Value = "RPR009"
$ws.Range("Y2").Value = "1234567RPR009"
$ws.Range("Z2").Value = "1234567RPR009"

$ws.Range("A2").Select()
